# Engage_NSReadings.xlsx — "SoGo NS reading issues and DMX DP issues are fixed"
#
# Performance_TC1 (row 2 of the PerformanceTC sheet) is switched off:
#   - CaseToRun (B2): Y -> N
#   - Pass/Fail/Skip (J2): PASS -> SKIP (recolored to match the existing SKIP rows)
#   - TimeLoad (K2): 8.81 -> 0
# The other SKIP rows (J3:J12) are re-touched so their fill/font is refreshed
# to the same yellow "SKIP" styling, and the active selection moves to B2.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(4)   # "PerformanceTC"
$ws.Activate()

# --- CaseToRun flag for Performance_TC1 ---
$ws.Range("B2").Value = "N"

# --- Pass/Fail/Skip result for Performance_TC1 ---
$ws.Range("J2").Value = "SKIP"

# Recolor J2:J12 with the same "SKIP" look (black text on yellow fill) that
# the rest of the column already used, so the freshly-SKIPped row 2 matches.
$skipRange = $ws.Range("J2:J12")
$skipRange.Font.Name = "Calibri"
$skipRange.Font.Size = 11
$skipRange.Font.Color = 0        # black
$skipRange.Interior.Color = 65535 # yellow

# --- TimeLoad for Performance_TC1: copy the existing textual "0" used by
#     the other rows so the cell stays a text value (not a numeric 0) ---
$ws.Range("K3").Copy()
$ws.Range("K2").PasteSpecial(-4163)  # xlPasteValues
$excel.CutCopyMode = $false

# --- move the active selection to B2 ---
$ws.Range("B2").Select()
